$d = $word.ActiveDocument

# Locate and remove the paragraph containing "Jkgklhjk" entirely
# (including its paragraph mark), reverting it back out of the document.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Jkgklhjk*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
